# Append the new IP record (row 5) to Sheet1, matching the existing rows'
# layout, values, formatting and hyperlink style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A / B: Sno + IP (kept in sync, idempotent if already present) ---
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "120.221.212.160"

# --- Column C: Status ---
$ws.Range("C5").Value = "Malicious"

# --- Column D: Link (plain text shown, with a real hyperlink relationship) ---
$ws.Range("D5").Value = "https://www.virustotal.com/gui/ip-address/120.221.212.160/detection"
$ws.Hyperlinks.Add($ws.Range("D5"), "https://www.virustotal.com/gui/ip-address/120.221.212.160/detection")
# Re-apply the same "Hyperlink" cell style used by the other rows in this
# column (Hyperlinks.Add mints its own style variant by default).
$ws.Range("D5").Style = $ws.Range("D4").Style

# --- Column E: last_analysis_stats ---
$ws.Range("E5").Value = "{'harmless': 56, 'malicious': 12, 'suspicious': 1, 'undetected': 21, 'timeout': 0}"

# --- Column F: Country ---
$ws.Range("F5").Value = "China"

# --- Column G / H: whois_date / Last_Modification_Date (serial datetimes) ---
$ws.Range("G5").Value = 45328.57457175926
$ws.Range("H5").Value = 45340.61069444445
$ws.Range("G5:H5").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# --- Column I: AS_Owner ---
$ws.Range("I5").Value = "China Mobile Communications Group Co., Ltd."
